$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate (by paragraph index, from the bottom, skipping the very
# first / title paragraph) the bold "Play Cosmic Jewels Free and
# Enjoy Space Adventure" paragraph near the end of the document.
# ------------------------------------------------------------------

$totalParas = $d.Paragraphs.Count
$boldParaIndex = -1
for ($i = $totalParas; $i -ge 1; $i--) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.StartsWith("Play Cosmic Jewels Free and Enjoy Space Adventure") -and $i -ne 1) {
        $boldParaIndex = $i
        break
    }
}

# ------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document. We clone the
#    formatting/structure (leading empty run + bold run) of the
#    paragraph located above by copying its FormattedText, then
#    patch up the style and the actual wording.
# ------------------------------------------------------------------

$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Style = $d.Styles.Item("Normal")

# paragraph indices below this point have shifted down by one because
# of the paragraph we just inserted
$boldPara = $d.Paragraphs.Item($boldParaIndex + 1)
$metaPara.Range.FormattedText = $boldPara.Range.FormattedText

$oldBoldText = "Play Cosmic Jewels Free and Enjoy Space Adventure"
$newBoldText = "Meta description"
$boldRunRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start + $oldBoldText.Length)
$boldRunRange.Text = $newBoldText

$afterBoldPos = $metaPara.Range.Start + $newBoldText.Length
$restRange = $d.Range($afterBoldPos, $afterBoldPos)
$restRange.InsertAfter(": Read our review of Cosmic Jewels and play for free to explore space, win cash and trigger special features.")

# ------------------------------------------------------------------
# 2. Near the end of the document, remove the (now duplicated) bold
#    "Play Cosmic Jewels Free and Enjoy Space Adventure" paragraph
#    entirely, and replace the text of the following italic
#    paragraph with the new image-prompt copy.
# ------------------------------------------------------------------

$searchRange = $d.Content
$found = $true
$lastMatch = $null
while ($found) {
    $found = $searchRange.Find.Execute("Play Cosmic Jewels Free and Enjoy Space Adventure")
    if ($found) {
        $lastMatch = $d.Range($searchRange.Start, $searchRange.End)
        $searchRange.Collapse(0)
    }
}

$dupPara = $lastMatch.Paragraphs.Item(1)
$dupParaFull = $d.Range($dupPara.Range.Start, $dupPara.Range.End + 1)
$dupParaFull.Delete()

$oldText = "Read our review of Cosmic Jewels and play for free to explore space, win cash and trigger special features."
$newText = "Create an eye-catching feature image for Cosmic Jewels that captures the excitement and adventure of the game. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be holding a space helmet in one hand and a shining diamond in the other, with stars and planets in the background. Use bold colors and dynamic poses to make the image stand out and draw players in. Make sure the image accurately represents the space adventure theme of the game and entices players to join in the hunt for valuable jewels."

$searchRange2 = $d.Content
$found2 = $true
$lastOldMatch = $null
while ($found2) {
    $found2 = $searchRange2.Find.Execute($oldText)
    if ($found2) {
        $lastOldMatch = $d.Range($searchRange2.Start, $searchRange2.End)
        $searchRange2.Collapse(0)
    }
}
$lastOldMatch.Text = $newText
